$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Hours" values for the two log entries
$ws.Range("C3").Value = "9 to 5"
$ws.Range("C3").NumberFormat = "d-mmm"
$ws.Range("C2").Value = "9 to 6"

# Remove the hyperlinks from the References column (displayed text stays, link removed)
$ws.Hyperlinks.Delete()
